# Refresh the NATMI ligand-receptor pair table (Bmp2-Bmpr1b) with the recomputed TPM-based
# values.  This reproduces the new 4x3 sending/target-cluster cross-product (ECs, FAPs,
# MuSCs, Resolving-Mac sending to FAPs, MuSCs, Resolving-Mac) together with the refreshed
# expression/specificity metrics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Bmp2"
$ws.Cells.Item(2, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(2, 4).Value2 = "FAPs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 1.729584333333333
$ws.Cells.Item(2, 8).Value2 = 5.188753
$ws.Cells.Item(2, 9).Value2 = 0.2476387648475193
$ws.Cells.Item(2, 10).Value2 = 0.2476387648475193
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 3.510190333333334
$ws.Cells.Item(2, 14).Value2 = 10.530571
$ws.Cells.Item(2, 15).Value2 = 0.7927501397588634
$ws.Cells.Item(2, 16).Value2 = 0.7927501397588635
$ws.Cells.Item(2, 17).Value2 = 6.071170207551445
$ws.Cells.Item(2, 18).Value2 = 54.64053186796301
$ws.Cells.Item(2, 19).Value2 = 0.1963156654425832
$ws.Cells.Item(2, 20).Value2 = 0.1963156654425832

# Row 3: ECs -> MuSCs
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Bmp2"
$ws.Cells.Item(3, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(3, 4).Value2 = "MuSCs"
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 1.729584333333333
$ws.Cells.Item(3, 8).Value2 = 5.188753
$ws.Cells.Item(3, 9).Value2 = 0.2476387648475193
$ws.Cells.Item(3, 10).Value2 = 0.2476387648475193
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 0.8869683333333332
$ws.Cells.Item(3, 14).Value2 = 2.660905
$ws.Cells.Item(3, 15).Value2 = 0.2003151406163121
$ws.Cells.Item(3, 16).Value2 = 0.2003151406163121
$ws.Cells.Item(3, 17).Value2 = 1.534086533496111
$ws.Cells.Item(3, 18).Value2 = 13.806778801465
$ws.Cells.Item(3, 19).Value2 = 0.04960579400248067
$ws.Cells.Item(3, 20).Value2 = 0.04960579400248068

# Row 4: ECs -> Resolving-Mac
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Bmp2"
$ws.Cells.Item(4, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(4, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 1.729584333333333
$ws.Cells.Item(4, 8).Value2 = 5.188753
$ws.Cells.Item(4, 9).Value2 = 0.2476387648475193
$ws.Cells.Item(4, 10).Value2 = 0.2476387648475193
$ws.Cells.Item(4, 11).Value2 = 1
$ws.Cells.Item(4, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(4, 13).Value2 = 0.030706
$ws.Cells.Item(4, 14).Value2 = 0.09211800000000001
$ws.Cells.Item(4, 15).Value2 = 0.006934719624824425
$ws.Cells.Item(4, 16).Value2 = 0.006934719624824427
$ws.Cells.Item(4, 17).Value2 = 0.05310861653933334
$ws.Cells.Item(4, 18).Value2 = 0.4779775488540001
$ws.Cells.Item(4, 19).Value2 = 0.001717305402455373
$ws.Cells.Item(4, 20).Value2 = 0.001717305402455373

# Row 5: FAPs -> FAPs
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Bmp2"
$ws.Cells.Item(5, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(5, 4).Value2 = "FAPs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 1.611936333333333
$ws.Cells.Item(5, 8).Value2 = 4.835809
$ws.Cells.Item(5, 9).Value2 = 0.2307941364328804
$ws.Cells.Item(5, 10).Value2 = 0.2307941364328804
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 3.510190333333334
$ws.Cells.Item(5, 14).Value2 = 10.530571
$ws.Cells.Item(5, 15).Value2 = 0.7927501397588634
$ws.Cells.Item(5, 16).Value2 = 0.7927501397588635
$ws.Cells.Item(5, 17).Value2 = 5.658203335215445
$ws.Cells.Item(5, 18).Value2 = 50.923830016939
$ws.Cells.Item(5, 19).Value2 = 0.1829620839126921
$ws.Cells.Item(5, 20).Value2 = 0.1829620839126921

# Row 6: FAPs -> MuSCs
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Bmp2"
$ws.Cells.Item(6, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(6, 4).Value2 = "MuSCs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 1.611936333333333
$ws.Cells.Item(6, 8).Value2 = 4.835809
$ws.Cells.Item(6, 9).Value2 = 0.2307941364328804
$ws.Cells.Item(6, 10).Value2 = 0.2307941364328804
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 0.8869683333333332
$ws.Cells.Item(6, 14).Value2 = 2.660905
$ws.Cells.Item(6, 15).Value2 = 0.2003151406163121
$ws.Cells.Item(6, 16).Value2 = 0.2003151406163121
$ws.Cells.Item(6, 17).Value2 = 1.429736483016111
$ws.Cells.Item(6, 18).Value2 = 12.867628347145
$ws.Cells.Item(6, 19).Value2 = 0.04623155989297275
$ws.Cells.Item(6, 20).Value2 = 0.04623155989297276

# Row 7: FAPs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Bmp2"
$ws.Cells.Item(7, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(7, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 1.611936333333333
$ws.Cells.Item(7, 8).Value2 = 4.835809
$ws.Cells.Item(7, 9).Value2 = 0.2307941364328804
$ws.Cells.Item(7, 10).Value2 = 0.2307941364328804
$ws.Cells.Item(7, 11).Value2 = 1
$ws.Cells.Item(7, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(7, 13).Value2 = 0.030706
$ws.Cells.Item(7, 14).Value2 = 0.09211800000000001
$ws.Cells.Item(7, 15).Value2 = 0.006934719624824425
$ws.Cells.Item(7, 16).Value2 = 0.006934719624824427
$ws.Cells.Item(7, 17).Value2 = 0.04949611705133334
$ws.Cells.Item(7, 18).Value2 = 0.445465053462
$ws.Cells.Item(7, 19).Value2 = 0.001600492627215502
$ws.Cells.Item(7, 20).Value2 = 0.001600492627215502

# Row 8: MuSCs -> FAPs
$ws.Cells.Item(8, 1).Value2 = "MuSCs"
$ws.Cells.Item(8, 2).Value2 = "Bmp2"
$ws.Cells.Item(8, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(8, 4).Value2 = "FAPs"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 2.743651333333334
$ws.Cells.Item(8, 8).Value2 = 8.230954000000001
$ws.Cells.Item(8, 9).Value2 = 0.3928310486309039
$ws.Cells.Item(8, 10).Value2 = 0.3928310486309038
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 3.510190333333334
$ws.Cells.Item(8, 14).Value2 = 10.530571
$ws.Cells.Item(8, 15).Value2 = 0.7927501397588634
$ws.Cells.Item(8, 16).Value2 = 0.7927501397588635
$ws.Cells.Item(8, 17).Value2 = 9.630738388303779
$ws.Cells.Item(8, 18).Value2 = 86.676645494734
$ws.Cells.Item(8, 19).Value2 = 0.3114168687037699
$ws.Cells.Item(8, 20).Value2 = 0.3114168687037699

# Row 9: MuSCs -> MuSCs
$ws.Cells.Item(9, 1).Value2 = "MuSCs"
$ws.Cells.Item(9, 2).Value2 = "Bmp2"
$ws.Cells.Item(9, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(9, 4).Value2 = "MuSCs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 2.743651333333334
$ws.Cells.Item(9, 8).Value2 = 8.230954000000001
$ws.Cells.Item(9, 9).Value2 = 0.3928310486309039
$ws.Cells.Item(9, 10).Value2 = 0.3928310486309038
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 0.8869683333333332
$ws.Cells.Item(9, 14).Value2 = 2.660905
$ws.Cells.Item(9, 15).Value2 = 0.2003151406163121
$ws.Cells.Item(9, 16).Value2 = 0.2003151406163121
$ws.Cells.Item(9, 17).Value2 = 2.433531850374445
$ws.Cells.Item(9, 18).Value2 = 21.90178665337
$ws.Cells.Item(9, 19).Value2 = 0.07869000674495284
$ws.Cells.Item(9, 20).Value2 = 0.07869000674495284

# Row 10: MuSCs -> Resolving-Mac
$ws.Cells.Item(10, 1).Value2 = "MuSCs"
$ws.Cells.Item(10, 2).Value2 = "Bmp2"
$ws.Cells.Item(10, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(10, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 2.743651333333334
$ws.Cells.Item(10, 8).Value2 = 8.230954000000001
$ws.Cells.Item(10, 9).Value2 = 0.3928310486309039
$ws.Cells.Item(10, 10).Value2 = 0.3928310486309038
$ws.Cells.Item(10, 11).Value2 = 1
$ws.Cells.Item(10, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(10, 13).Value2 = 0.030706
$ws.Cells.Item(10, 14).Value2 = 0.09211800000000001
$ws.Cells.Item(10, 15).Value2 = 0.006934719624824425
$ws.Cells.Item(10, 16).Value2 = 0.006934719624824427
$ws.Cells.Item(10, 17).Value2 = 0.08424655784133335
$ws.Cells.Item(10, 18).Value2 = 0.7582190205720001
$ws.Cells.Item(10, 19).Value2 = 0.002724173182181087
$ws.Cells.Item(10, 20).Value2 = 0.002724173182181088

# Row 11: Resolving-Mac -> FAPs
$ws.Cells.Item(11, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value2 = "Bmp2"
$ws.Cells.Item(11, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(11, 4).Value2 = "FAPs"
$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 6).Value2 = 1
$ws.Cells.Item(11, 7).Value2 = 0.8991316666666668
$ws.Cells.Item(11, 8).Value2 = 2.697395
$ws.Cells.Item(11, 9).Value2 = 0.1287360500886965
$ws.Cells.Item(11, 10).Value2 = 0.1287360500886965
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 12).Value2 = 1
$ws.Cells.Item(11, 13).Value2 = 3.510190333333334
$ws.Cells.Item(11, 14).Value2 = 10.530571
$ws.Cells.Item(11, 15).Value2 = 0.7927501397588634
$ws.Cells.Item(11, 16).Value2 = 0.7927501397588635
$ws.Cells.Item(11, 17).Value2 = 3.156123284727223
$ws.Cells.Item(11, 18).Value2 = 28.405109562545
$ws.Cells.Item(11, 19).Value2 = 0.1020555216998182
$ws.Cells.Item(11, 20).Value2 = 0.1020555216998182

# Row 12: Resolving-Mac -> MuSCs
$ws.Cells.Item(12, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value2 = "Bmp2"
$ws.Cells.Item(12, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(12, 4).Value2 = "MuSCs"
$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 6).Value2 = 1
$ws.Cells.Item(12, 7).Value2 = 0.8991316666666668
$ws.Cells.Item(12, 8).Value2 = 2.697395
$ws.Cells.Item(12, 9).Value2 = 0.1287360500886965
$ws.Cells.Item(12, 10).Value2 = 0.1287360500886965
$ws.Cells.Item(12, 11).Value2 = 3
$ws.Cells.Item(12, 12).Value2 = 1
$ws.Cells.Item(12, 13).Value2 = 0.8869683333333332
$ws.Cells.Item(12, 14).Value2 = 2.660905
$ws.Cells.Item(12, 15).Value2 = 0.2003151406163121
$ws.Cells.Item(12, 16).Value2 = 0.2003151406163121
$ws.Cells.Item(12, 17).Value2 = 0.7975013158305556
$ws.Cells.Item(12, 18).Value2 = 7.177511842475
$ws.Cells.Item(12, 19).Value2 = 0.02578777997590584
$ws.Cells.Item(12, 20).Value2 = 0.02578777997590584

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(13, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value2 = "Bmp2"
$ws.Cells.Item(13, 3).Value2 = "Bmpr1b"
$ws.Cells.Item(13, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 6).Value2 = 1
$ws.Cells.Item(13, 7).Value2 = 0.8991316666666668
$ws.Cells.Item(13, 8).Value2 = 2.697395
$ws.Cells.Item(13, 9).Value2 = 0.1287360500886965
$ws.Cells.Item(13, 10).Value2 = 0.1287360500886965
$ws.Cells.Item(13, 11).Value2 = 1
$ws.Cells.Item(13, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(13, 13).Value2 = 0.030706
$ws.Cells.Item(13, 14).Value2 = 0.09211800000000001
$ws.Cells.Item(13, 15).Value2 = 0.006934719624824425
$ws.Cells.Item(13, 16).Value2 = 0.006934719624824427
$ws.Cells.Item(13, 17).Value2 = 0.02760873695666667
$ws.Cells.Item(13, 18).Value2 = 0.24847863261
$ws.Cells.Item(13, 19).Value2 = 0.0008927484129724639
$ws.Cells.Item(13, 20).Value2 = 0.0008927484129724641
